$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "Follow up questions" sheet: insert new follow-up bullet questions
# and re-flow the existing content down to make room for them.
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Follow up questions")

# Make room for the three new bullet rows (insert from the top down so
# each subsequent row number below still lines up with its target row).
$ws4.Rows("8:8").Insert() | Out-Null
$ws4.Rows("12:12").Insert() | Out-Null
$ws4.Rows("16:16").Insert() | Out-Null

# New bullet appended after "4) existing process/system assessment".
$ws4.Range("C16").Value = " - what is the business definition of each attributes?"

# New bullet under "3) stakeholder assessment".
$ws4.Range("C12").Value = " - why is this important?"

# New bullet under "2) feature assesment".
$ws4.Range("C8").Value = " - is there any alternative feature other than bulk upload that could lead to the same objective? (e.g. RFID automatic tracking system)"

$ws4.Activate() | Out-Null
$ws4.Range("C9").Select() | Out-Null

# ------------------------------------------------------------------
# "User Story" sheet: selection moved.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("User Story")
$ws1.Activate() | Out-Null
$ws1.Range("E23").Select() | Out-Null
